$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.168.10"
$ws.Range("E2").Value = "'  -1.55%  "
$ws.Range("D3").Value = "'2.357.18"
$ws.Range("E3").Value = "'  +1.92%  "
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'301.96"
$ws.Range("E5").Value = "'  +1.09%  "
$ws.Range("D6").Value = "'99.78"
$ws.Range("E6").Value = "'  +0.68%  "
$ws.Range("E7").Value = "'  -0.25%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E9").Value = "'  -3.02%  "
$ws.Range("D10").Value = "'34.54"
$ws.Range("E10").Value = "'  -3.41%  "
$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "'  +0.14%  "
$ws.Range("D12").Value = "'7.15"
$ws.Range("E12").Value = "'  -2.67%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "'  -0.25%  "
$ws.Range("D14").Value = "'2.717.20"
$ws.Range("E14").Value = "'  +1.96%  "
$ws.Range("D15").Value = "'2.358.84"
$ws.Range("E15").Value = "'  +2.10%  "
$ws.Range("D16").Value = "'0.811"
$ws.Range("E16").Value = "'  -0.58%  "
$ws.Range("D17").Value = "'13.61"
$ws.Range("E17").Value = "'  -2.78%  "
$ws.Range("D18").Value = "'46.093.48"
$ws.Range("E18").Value = "'  -1.55%  "
$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "'  -2.71%  "
$ws.Range("D20").Value = "'0.0₃0966"
$ws.Range("E20").Value = "'  +2.81%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("E21").Value = "'  -0.76%  "
$ws.Range("D22").Value = "'67.48"
$ws.Range("E22").Value = "'  +1.03%  "
$ws.Range("D23").Value = "'246.67"
$ws.Range("E23").Value = "'  -0.91%  "
$ws.Range("D24").Value = "'2.84"
$ws.Range("E24").Value = "'  -2.67%  "
$ws.Range("E25").Value = "'  -0.25%  "
$ws.Range("E26").Value = "'  -3.05%  "
$ws.Range("D27").Value = "'39.76"
$ws.Range("E27").Value = "'  -6.65%  "
$ws.Range("E28").Value = "'  -2.80%  "
$ws.Range("D29").Value = "'9.81"
$ws.Range("E29").Value = "'  -0.43%  "
$ws.Range("D30").Value = "'21.03"
$ws.Range("E30").Value = "'  +4.12%  "
$ws.Range("D31").Value = "'3.76"
$ws.Range("E31").Value = "'  +20.75%  "
$ws.Range("E32").Value = "'  +6.19%  "
$ws.Range("D33").Value = "'5.54"
$ws.Range("E33").Value = "'  -3.51%  "
$ws.Range("D34").Value = "'146.29"
$ws.Range("E34").Value = "'  -0.67%  "
$ws.Range("D35").Value = "'0.0777"
$ws.Range("E35").Value = "'  -2.64%  "
$ws.Range("E36").Value = "'  -1.15%  "
$ws.Range("D37").Value = "'1.90"
$ws.Range("E37").Value = "'  +5.45%  "
$ws.Range("E38").Value = "'  -1.88%  "
$ws.Range("D39").Value = "'15.08"
$ws.Range("E39").Value = "'  -4.42%  "
$ws.Range("D40").Value = "'3.96"
$ws.Range("E40").Value = "'  -1.12%  "
$ws.Range("D41").Value = "'0.0300"
$ws.Range("E41").Value = "'  -2.32%  "
$ws.Range("D42").Value = "'3.24"
$ws.Range("E42").Value = "'  -5.09%  "
$ws.Range("D43").Value = "'1.896.05"
$ws.Range("E43").Value = "'  +2.85%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Range("D45").Value = "'92.51"
$ws.Range("E45").Value = "'  +1.47%  "
$ws.Range("D46").Value = "'1.81"
$ws.Range("E46").Value = "'  -8.76%  "
$ws.Range("D47").Value = "'0.187"
$ws.Range("E47").Value = "'  -6.82%  "
$ws.Range("E48").Value = "'  +3.43%  "
$ws.Range("D49").Value = "'97.78"
$ws.Range("E49").Value = "'  +0.53%  "
$ws.Range("D50").Value = "'2.588.48"
$ws.Range("E50").Value = "'  +1.71%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'69.15"
$ws.Range("E51").Value = "'  -8.76%  "
